$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before column F (Nombre del material), which will
# become the "Precio Residencial", "Precio Comercial" and "Precio Distribuidor"
# price-tier columns.
$ws.Range("F1:H1").EntireColumn.Insert()

# New header cells (row 2)
$ws.Range("F2").Value = "Precio Residencial"
$ws.Range("G2").Value = "Precio Comercial"
$ws.Range("H2").Value = "Precio Distribuidor"

# New data cells (row 3)
$ws.Range("F3").Value = 250
$ws.Range("G3").Value = 245
$ws.Range("H3").Value = 240

# Size the new price columns to fit their header text (matches Excel's
# "best fit" auto-sizing for the inserted columns)
$ws.Range("F1").EntireColumn.ColumnWidth = 16.428571428571427
$ws.Range("G1").EntireColumn.ColumnWidth = 15.0
$ws.Range("H1").EntireColumn.ColumnWidth = 16.857142857142858
